$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 <- original row 6 data
$ws.Range("A5").Value = 130963816
$ws.Range("B5").Value = 79243
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("M5").ClearContents()
$ws.Range("Q5").Value = 445932
$ws.Range("R5").Value = 6760079
$ws.Range("Z5").Value = "14:08"
$ws.Range("AB5").Value = "14:08"
$ws.Range("AC5").Value = "Rikligt i närområdet"

# Row 6 <- original row 5 data
$ws.Range("A6").Value = 130960607
$ws.Range("B6").Value = 57884
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("M6").Value = "äldre spår"
$ws.Range("Q6").Value = 446240
$ws.Range("R6").Value = 6759818
$ws.Range("Z6").Value = "10:26"
$ws.Range("AB6").Value = "10:26"
$ws.Range("AC6").ClearContents()

# Row 9 <- original row 10 data
$ws.Range("A9").Value = 130962883
$ws.Range("B9").Value = 79243
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("M9").ClearContents()
$ws.Range("Q9").Value = 445987
$ws.Range("R9").Value = 6759938
$ws.Range("Z9").Value = "10:26"
$ws.Range("AB9").Value = "10:26"
$ws.Range("AC9").ClearContents()

# Row 10 <- original row 9 data
$ws.Range("A10").Value = 130960395
$ws.Range("B10").Value = 8451
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 106545
$ws.Range("F10").Value = "Mindre märgborre"
$ws.Range("G10").Value = "Tomicus minor"
$ws.Range("H10").Value = "(Hartig, 1834)"
$ws.Range("M10").Value = "äldre gnagspår"
$ws.Range("Q10").Value = 446272
$ws.Range("R10").Value = 6759739
$ws.Range("Z10").Value = "10:26"
$ws.Range("AB10").Value = "10:26"
$ws.Range("AC10").ClearContents()

# Row 11 <- original row 12 data
$ws.Range("A11").Value = 130961179
$ws.Range("B11").Value = 79862
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6453
$ws.Range("F11").Value = "Vedskivlav"
$ws.Range("G11").Value = "Hertelidea botryosa"
$ws.Range("H11").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("M11").ClearContents()
$ws.Range("Q11").Value = 446122
$ws.Range("R11").Value = 6760020
$ws.Range("Z11").Value = "10:26"
$ws.Range("AB11").Value = "10:26"
$ws.Range("AC11").Value = "Ringhack på stam i bakgrund"

# Row 12 <- original row 11 data
$ws.Range("A12").Value = 130961218
$ws.Range("B12").Value = 57884
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = "Tretåig hackspett"
$ws.Range("G12").Value = "Picoides tridactylus"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("M12").Value = "äldre spår"
$ws.Range("Q12").Value = 446122
$ws.Range("R12").Value = 6760020
$ws.Range("Z12").Value = "10:26"
$ws.Range("AB12").Value = "10:26"
$ws.Range("AC12").ClearContents()

# Row 18 <- original row 19 data
$ws.Range("A18").Value = 130961956
$ws.Range("B18").Value = 79862
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 6453
$ws.Range("F18").Value = "Vedskivlav"
$ws.Range("G18").Value = "Hertelidea botryosa"
$ws.Range("H18").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("M18").ClearContents()
$ws.Range("Q18").Value = 446084
$ws.Range("R18").Value = 6759981
$ws.Range("Z18").Value = "10:26"
$ws.Range("AB18").Value = "10:26"
$ws.Range("AC18").Value = "Miljöbilder"

# Row 19 <- original row 18 data
$ws.Range("A19").Value = 130960843
$ws.Range("B19").Value = 79243
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6425
$ws.Range("F19").Value = "Garnlav"
$ws.Range("G19").Value = "Alectoria sarmentosa"
$ws.Range("H19").Value = "(Ach.) Ach."
$ws.Range("M19").ClearContents()
$ws.Range("Q19").Value = 446247
$ws.Range("R19").Value = 6759903
$ws.Range("Z19").Value = "10:26"
$ws.Range("AB19").Value = "10:26"
$ws.Range("AC19").ClearContents()

# Row 22 <- original row 23 data
$ws.Range("A22").Value = 130962722
$ws.Range("B22").Value = 79862
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6453
$ws.Range("F22").Value = "Vedskivlav"
$ws.Range("G22").Value = "Hertelidea botryosa"
$ws.Range("H22").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("M22").ClearContents()
$ws.Range("Q22").Value = 446008
$ws.Range("R22").Value = 6759948
$ws.Range("Z22").Value = "10:26"
$ws.Range("AB22").Value = "10:26"
$ws.Range("AC22").ClearContents()

# Row 23 <- original row 22 data
$ws.Range("A23").Value = 130963976
$ws.Range("B23").Value = 79243
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 6425
$ws.Range("F23").Value = "Garnlav"
$ws.Range("G23").Value = "Alectoria sarmentosa"
$ws.Range("H23").Value = "(Ach.) Ach."
$ws.Range("M23").ClearContents()
$ws.Range("Q23").Value = 445929
$ws.Range("R23").Value = 6760099
$ws.Range("Z23").Value = "14:08"
$ws.Range("AB23").Value = "14:08"
$ws.Range("AC23").Value = "Miljöbild"

# Row 27 <- original row 30 data
$ws.Range("A27").Value = 130962736
$ws.Range("B27").Value = 79833
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 229821
$ws.Range("F27").Value = "Vedflamlav"
$ws.Range("G27").Value = "Ramboldia elabens"
$ws.Range("H27").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("M27").ClearContents()
$ws.Range("Q27").Value = 446008
$ws.Range("R27").Value = 6759948
$ws.Range("Z27").Value = "10:26"
$ws.Range("AB27").Value = "10:26"
$ws.Range("AC27").ClearContents()

# Row 29 <- original row 27 data
$ws.Range("A29").Value = 130961750
$ws.Range("B29").Value = 79243
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("M29").ClearContents()
$ws.Range("Q29").Value = 446098
$ws.Range("R29").Value = 6760061
$ws.Range("Z29").Value = "10:26"
$ws.Range("AB29").Value = "10:26"
$ws.Range("AC29").Value = "Rikligt i en radie av ca 50 meter"

# Row 30 <- original row 29 data
$ws.Range("A30").Value = 130963807
$ws.Range("B30").Value = 57881
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 100049
$ws.Range("F30").Value = "Spillkråka"
$ws.Range("G30").Value = "Dryocopus martius"
$ws.Range("H30").Value = "(Linnaeus, 1758)"
$ws.Range("M30").Value = "färska spår"
$ws.Range("Q30").Value = 445932
$ws.Range("R30").Value = 6760079
$ws.Range("Z30").Value = "14:08"
$ws.Range("AB30").Value = "14:08"
$ws.Range("AC30").ClearContents()
